$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B17 (Compositional row) already holds the literal text "true" as a shared
# string. Copy it into B7 (Experimental) and B14 (Case Sensitive) so the new
# cells pick up the same shared-string text entry instead of Excel's
# automatic "true"/"false" -> Boolean literal coercion that a plain
# .Value assignment would trigger.
$ws.Range("B17").Copy()
$ws.Range("B7").PasteSpecial(-4163)
$ws.Range("B14").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# Update the Date value (CodeSystem "Date" property).
$ws.Range("B8").Value = "2024-02-19T18:37:26-06:00"
